$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to text so numeric-looking / percent strings are not
# auto-converted to numbers or dates by Excel when assigned via .Value
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "25.977.04"
$ws.Range("E2").Value = "  -2.38%  "
$ws.Range("D3").Value = "1.663.92"
$ws.Range("E3").Value = "  -1.98%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "215.81"
$ws.Range("E5").Value = "  -1.88%  "
$ws.Range("D6").Value = "0.5085"
$ws.Range("E6").Value = "  -1.46%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.2632"
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("D9").Value = "0.06374"
$ws.Range("E9").Value = "  +1.51%  "
$ws.Range("D10").Value = "21.79"
$ws.Range("E10").Value = "  -2.03%  "
$ws.Range("D11").Value = "0.07412"
$ws.Range("E11").Value = "  +0.69%  "
$ws.Range("D12").Value = "1.672.14"
$ws.Range("E12").Value = "  -1.61%  "
$ws.Range("D13").Value = "4.488"
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("D14").Value = "0.5805"
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("D15").Value = "0.000008441"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").Value = "63.99"
$ws.Range("E16").Value = "  -2.55%  "
$ws.Range("D17").Value = "26.029.71"
$ws.Range("E17").Value = "  -2.38%  "
$ws.Range("D18").Value = "4.905"
$ws.Range("E18").Value = "  -2.66%  "
$ws.Range("D19").Value = "1.006"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").Value = "10.66"
$ws.Range("E20").Value = "  -2.41%  "
$ws.Range("D21").Value = "188.73"
$ws.Range("E21").Value = "  +0.90%  "
$ws.Range("D22").Value = "6.183"
$ws.Range("E22").Value = "  -1.65%  "
$ws.Range("D23").Value = "1.006"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "144.80"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").Value = "7.579"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("D26").Value = "0.1187"
$ws.Range("E26").Value = "  +2.76%  "
$ws.Range("E27").Value = "  -1.56%  "
$ws.Range("D28").Value = "0.06575"
$ws.Range("E28").Value = "  +15.31%  "
$ws.Range("D29").Value = "1.311"
$ws.Range("E29").Value = "  -0.86%  "
$ws.Range("D30").Value = "1.313"
$ws.Range("E30").Value = "  -2.21%  "
$ws.Range("D31").Value = "3.521"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").Value = "3.498"
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("D33").Value = "1.626"
$ws.Range("E33").Value = "  -2.02%  "
$ws.Range("D34").Value = "1.013"
$ws.Range("E34").Value = "  -1.33%  "
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("D37").Value = "2.699"
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("D38").Value = "6.190"
$ws.Range("E38").Value = "  +5.41%  "
$ws.Range("D39").Value = "0.01603"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").Value = "1.072.16"
$ws.Range("E40").Value = "  -2.78%  "
$ws.Range("D41").Value = "0.8583"
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("E42").Value = "  +0.49%  "
$ws.Range("D43").Value = "100.21"
$ws.Range("E43").Value = "  +1.16%  "
$ws.Range("D44").Value = "1.811.54"
$ws.Range("E44").Value = "  -2.47%  "
$ws.Range("E45").Value = "  +4.22%  "
$ws.Range("D46").Value = "56.04"
$ws.Range("E46").Value = "  -1.39%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "7.994"
$ws.Range("E48").Value = "  -2.39%  "
$ws.Range("D49").Value = "0.05205"
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("D50").Value = "0.4290"
$ws.Range("E50").Value = "  -0.82%  "
$ws.Range("D51").Value = "5.933"
$ws.Range("E51").Value = "  +2.39%  "

# Restore default (no explicit) style on the touched range so cells
# revert to the original unstyled inline-string presentation
$ws.Range("D2:E51").Style = "Normal"
